# archivocuotas.xlsx - "Cuotas" sheet
# Commit: cambio logica de asignacion PREPAID (60) y mas de 1 cuota CREDIT (10).
#         Genero DEVBOTON vacio
#
# The new assignment logic produced a different set of transaction ids in
# column A (rows 2-7); column B ("Cuotas") is unchanged. The run also
# (re)generates an empty header/footer definition for the sheet ("DEVBOTON
# vacio").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated transaction ids (column A) -----------------------------------
$ws.Range("A2").Value = 1638041672.0
$ws.Range("A3").Value = 1611264484.0
$ws.Range("A4").Value = 1650221344.0
$ws.Range("A5").Value = 1908704539.0
$ws.Range("A6").Value = 1611905756.0
$ws.Range("A7").Value = 1772895896.0

# --- Regenerate an empty header/footer ("DEVBOTON vacio") -----------------
$ps = $ws.PageSetup
$ps.OddAndEvenPagesHeaderFooter = $false
$ps.DifferentFirstPageHeaderFooter = $false
$ps.ScaleWithDocHeaderFooter = $true
$ps.AlignMarginsHeaderFooter = $true
$ps.LeftHeader = ""
$ps.CenterHeader = ""
$ps.RightHeader = ""
$ps.LeftFooter = ""
$ps.CenterFooter = ""
$ps.RightFooter = ""
